$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text columns to Text format so leading zeros / numeric-looking IDs are preserved
$ws.Range("A2:E13").NumberFormat = "@"
$ws.Range("G2:I13").NumberFormat = "@"
$ws.Range("L2:M13").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "HECTOR ALVARADO"
$ws.Range("B2").Value = "0401199021421"
$ws.Range("C2").Value = "Ningun parentesco"
$ws.Range("D2").Value = "Oficina Principal"
$ws.Range("E2").Value = "005-001-000000052"
$ws.Range("F2").Value = 117
$ws.Range("G2").Value = "Herencia"
$ws.Range("H2").Value = "Esta en el trabajo/negocio"
$ws.Range("I2").Value = "Depósito"
$ws.Range("J2").Value = 25000
$ws.Range("K2").Value = 44720.596666666665
$ws.Range("L2").Value = "MJVG"
$ws.Range("M2").Value = ""

# Row 3
$ws.Range("A3").Value = "JUAN DAVID PEREZ MENDOZA "
$ws.Range("B3").Value = "0303030303033"
$ws.Range("C3").Value = "Sobrino (a)"
$ws.Range("D3").Value = "Oficina Principal"
$ws.Range("E3").Value = "005-001-000002819"
$ws.Range("F3").Value = 4671
$ws.Range("G3").Value = "Por Compra o Venta de Naranja"
$ws.Range("H3").Value = "Falta de tiempo"
$ws.Range("I3").Value = "Depósito"
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 44720.38976851852
$ws.Range("L3").Value = "MJVG"
$ws.Range("M3").Value = ""

# Row 4
$ws.Range("A4").Value = "JUAN JOSE PEREZ RAMIREZ "
$ws.Range("B4").Value = "0000100012345"
$ws.Range("C4").Value = "Hermano (a)"
$ws.Range("D4").Value = "Oficina Principal"
$ws.Range("E4").Value = "005-001-000002598"
$ws.Range("F4").Value = 4367
$ws.Range("G4").Value = "Venta de una propiedad"
$ws.Range("H4").Value = "Esta en el trabajo/negocio"
$ws.Range("I4").Value = "Operación de Unired"
$ws.Range("J4").Value = 17000
$ws.Range("K4").Value = 44714.602638888886
$ws.Range("L4").Value = "MJVG"
$ws.Range("M4").Value = "Prueba de ingresos menores con cliente nuevo "

# Row 5
$ws.Range("A5").Value = "MANU VILLEDA"
$ws.Range("B5").Value = "1415-1998-01225"
$ws.Range("C5").Value = "Sobrino (a)"
$ws.Range("D5").Value = "Oficina Principal"
$ws.Range("E5").Value = "005-1-2598"
$ws.Range("F5").Value = 4368
$ws.Range("G5").Value = "Préstamo obtenido"
$ws.Range("H5").Value = "No queria venir"
$ws.Range("I5").Value = "Depósito"
$ws.Range("J5").Value = 100000
$ws.Range("K5").Value = 44713.413310185184
$ws.Range("L5").Value = "MJVG"
$ws.Range("M5").Value = "prueba"

# Row 6
$ws.Range("A6").Value = "NICOLL MARTÍNEZ "
$ws.Range("B6").Value = "0401-2001-01227"
$ws.Range("C6").Value = "Hermano (a)"
$ws.Range("D6").Value = "Oficina Principal"
$ws.Range("E6").Value = "005-001-14759"
$ws.Range("F6").Value = 15
$ws.Range("G6").Value = "Venta de una propiedad"
$ws.Range("H6").Value = "Falta de tiempo"
$ws.Range("I6").Value = "Depósito"
$ws.Range("J6").Value = 150001
$ws.Range("K6").Value = 44711.483252314814
$ws.Range("L6").Value = "MJVG"
$ws.Range("M6").Value = "hola"

# Row 7
$ws.Range("A7").Value = "MANU VILLEDA"
$ws.Range("B7").Value = "1415-1998-01225"
$ws.Range("C7").Value = "Hermano (a)"
$ws.Range("D7").Value = "Oficina Principal"
$ws.Range("E7").Value = "005-002-14578"
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = "Venta de una propiedad"
$ws.Range("H7").Value = "Falta de tiempo"
$ws.Range("I7").Value = "Depósito"
$ws.Range("J7").Value = 145000
$ws.Range("K7").Value = 44711.45949074074
$ws.Range("L7").Value = "MJVG"
$ws.Range("M7").Value = "hola"

# Row 8
$ws.Range("A8").Value = "MANU VILLEDA"
$ws.Range("B8").Value = "1415-1998-01225"
$ws.Range("C8").Value = "Hermano (a)"
$ws.Range("D8").Value = "Oficina Principal"
$ws.Range("E8").Value = "005-002-17589"
$ws.Range("F8").Value = 15
$ws.Range("G8").Value = "Venta de una propiedad"
$ws.Range("H8").Value = "Falta de tiempo"
$ws.Range("I8").Value = "Depósito"
$ws.Range("J8").Value = 150000
$ws.Range("K8").Value = 44711.45756944444
$ws.Range("L8").Value = "MJVG"
$ws.Range("M8").Value = "Hola"

# Row 9
$ws.Range("A9").Value = "JESÚS GÓMEZ"
$ws.Range("B9").Value = "0401-1998-01225"
$ws.Range("C9").Value = "Hermano (a)"
$ws.Range("D9").Value = "Oficina Principal"
$ws.Range("E9").Value = "005-001-1759"
$ws.Range("F9").Value = 15
$ws.Range("G9").Value = "Venta de una propiedad"
$ws.Range("H9").Value = "Falta de tiempo"
$ws.Range("I9").Value = "Depósito"
$ws.Range("J9").Value = 13
$ws.Range("K9").Value = 44708.66537037037
$ws.Range("L9").Value = "MJVG"
$ws.Range("M9").Value = "Prueba 3"

# Row 10
$ws.Range("A10").Value = "MANU VILLEDA"
$ws.Range("B10").Value = "1415-1998-01225"
$ws.Range("C10").Value = "Hermano (a)"
$ws.Range("D10").Value = "Oficina Principal"
$ws.Range("E10").Value = "005-001-002"
$ws.Range("F10").Value = 15
$ws.Range("G10").Value = "Venta de una propiedad"
$ws.Range("H10").Value = "Falta de tiempo"
$ws.Range("I10").Value = "Depósito"
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = 44708.64915509259
$ws.Range("L10").Value = "MJVG"
$ws.Range("M10").Value = "Prueba de fecha"

# Row 11
$ws.Range("A11").Value = "MANU VILLEDA"
$ws.Range("B11").Value = "1415-1998-01225"
$ws.Range("C11").Value = "Hermano (a)"
$ws.Range("D11").Value = "Oficina Principal"
$ws.Range("E11").Value = "005-001-14151998"
$ws.Range("F11").Value = 15
$ws.Range("G11").Value = "Venta de una propiedad"
$ws.Range("H11").Value = "Falta de tiempo"
$ws.Range("I11").Value = "Depósito"
$ws.Range("J11").Value = 14500
$ws.Range("K11").Value = 44704.61928240741
$ws.Range("L11").Value = "MJVG"
$ws.Range("M11").Value = ""

# Row 12
$ws.Range("A12").Value = "MANU VILLEDA"
$ws.Range("B12").Value = "1415-1998-01225"
$ws.Range("C12").Value = "Hermano (a)"
$ws.Range("D12").Value = "Oficina Principal"
$ws.Range("E12").Value = "005-001-002"
$ws.Range("F12").Value = 15
$ws.Range("G12").Value = "Venta de una propiedad"
$ws.Range("H12").Value = "Falta de tiempo"
$ws.Range("I12").Value = "Depósito"
$ws.Range("J12").Value = 14
$ws.Range("K12").Value = 44704.373877314814
$ws.Range("L12").Value = "MJVG"
$ws.Range("M12").Value = "mamama"

# Row 13
$ws.Range("A13").Value = "MANU VILLEDA"
$ws.Range("B13").Value = "1415-1998-01225"
$ws.Range("C13").Value = "Hermano (a)"
$ws.Range("D13").Value = "Oficina Principal"
$ws.Range("E13").Value = "005-001-141597001"
$ws.Range("F13").Value = 15
$ws.Range("G13").Value = "Venta de una propiedad"
$ws.Range("H13").Value = "Falta de tiempo"
$ws.Range("I13").Value = "Depósito"
$ws.Range("J13").Value = 14500
$ws.Range("K13").Value = 44702.721979166665
$ws.Range("L13").Value = "MJVG"
$ws.Range("M13").Value = "Nada"
